# Carta Modelo Designacion Revisor - dynamic designation letter generation.
# Replace the fixed date / reviewer / thesis-title / student / tutor-email
# placeholders with the new values from the commit.
#
# We search with Find.Execute (no replacement arg, so Word just selects /
# collapses the range to the hit) and then assign the new text straight to
# the found Range's .Text property. Doing the substitution this way (rather
# than passing the replacement text into Find.Execute) keeps Word's
# "replace as you type" AutoFormat (smart quotes, etc.) out of the way, so
# straight quotes stay straight.

$d = $word.ActiveDocument

$dq = [char]34

# 1) Letter date.
$r = $d.Content
$r.Find.Execute("Puerto Ordaz, 03 noviembre 2022.") | Out-Null
$r.Text = "Puerto Ordaz, 29 de Julio del 2020"

# 2) Reviewer professor's name.
$r = $d.Content
$r.Find.Execute(" Franklin Bello") | Out-Null
$r.Text = " Larez, Jesus"

# 3) Thesis title.
$r = $d.Content
$r.Find.Execute($dq + "Creacion de robots para el mantenimiento de la biblioteca" + $dq + ". ") | Out-Null
$r.Text = $dq + "Titulo de propuesta de grado" + $dq + ". "

# 4) Student name mentioned next to the thesis title.
$r = $d.Content
$r.Find.Execute("Luis C Somoza; ") | Out-Null
$r.Text = "Somoza Ledezma, Luis Carlos"

# 5) Student name mentioned again further down (contact paragraph).
$r = $d.Content
$r.Find.Execute("Luis C Somoza") | Out-Null
$r.Text = "Somoza Ledezma, Luis Carlos"

# 6) Contact e-mail address.
$r = $d.Content
$r.Find.Execute("lmedinac@ucab.edu.ve ") | Out-Null
$r.Text = "lcmedina.19@ucab.edu.ve "
